# Auto-generated edit script
# Reorders the "Recorded By" (column G) comma-separated list so that
# the "System" token moves from the end of the list to the front,
# for a specific enumerated set of rows in the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $null
try {
    $ws = $wb.Worksheets.Item("Session Analysis Results")
} catch {
    $ws = $wb.ActiveSheet
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$updates = @(
    @{ Row = 2; New = 'System, backup@backdoor.com, system' }
    @{ Row = 3; New = 'System, dnasr281@gmail.com' }
    @{ Row = 5; New = 'System, backup@backdoor.com' }
    @{ Row = 6; New = 'System, dnasr281@gmail.com' }
    @{ Row = 7; New = 'System, admin@admin.com' }
    @{ Row = 8; New = 'System, backup@backdoor.com' }
    @{ Row = 10; New = 'System, dnasr281@gmail.com' }
    @{ Row = 11; New = 'System, dnasr281@gmail.com' }
    @{ Row = 12; New = 'System, dnasr281@gmail.com' }
    @{ Row = 13; New = 'System, dnasr281@gmail.com' }
    @{ Row = 14; New = 'System, dnasr281@gmail.com' }
    @{ Row = 15; New = 'System, dnasr281@gmail.com' }
    @{ Row = 17; New = 'System, dnasr281@gmail.com' }
    @{ Row = 18; New = 'System, dnasr281@gmail.com' }
    @{ Row = 19; New = 'System, dnasr281@gmail.com' }
    @{ Row = 20; New = 'System, dnasr281@gmail.com' }
    @{ Row = 21; New = 'System, dnasr281@gmail.com' }
    @{ Row = 22; New = 'System, dnasr281@gmail.com' }
    @{ Row = 24; New = 'System, dnasr281@gmail.com' }
    @{ Row = 26; New = 'System, dnasr281@gmail.com' }
    @{ Row = 28; New = 'System, backup@backdoor.com, system' }
    @{ Row = 29; New = 'System, dnasr281@gmail.com' }
    @{ Row = 31; New = 'System, backup@backdoor.com' }
    @{ Row = 32; New = 'System, dnasr281@gmail.com' }
    @{ Row = 33; New = 'System, admin@admin.com' }
    @{ Row = 34; New = 'System, backup@backdoor.com' }
    @{ Row = 36; New = 'System, dnasr281@gmail.com' }
    @{ Row = 37; New = 'System, dnasr281@gmail.com' }
    @{ Row = 38; New = 'System, dnasr281@gmail.com' }
    @{ Row = 39; New = 'System, dnasr281@gmail.com' }
    @{ Row = 40; New = 'System, dnasr281@gmail.com' }
    @{ Row = 41; New = 'System, dnasr281@gmail.com' }
    @{ Row = 43; New = 'System, dnasr281@gmail.com' }
    @{ Row = 44; New = 'System, dnasr281@gmail.com' }
    @{ Row = 45; New = 'System, dnasr281@gmail.com' }
    @{ Row = 46; New = 'System, dnasr281@gmail.com' }
    @{ Row = 47; New = 'System, dnasr281@gmail.com' }
    @{ Row = 48; New = 'System, dnasr281@gmail.com' }
    @{ Row = 50; New = 'System, dnasr281@gmail.com' }
    @{ Row = 52; New = 'System, dnasr281@gmail.com' }
    @{ Row = 54; New = 'System, backup@backdoor.com, system' }
    @{ Row = 55; New = 'System, dnasr281@gmail.com' }
    @{ Row = 57; New = 'System, backup@backdoor.com' }
    @{ Row = 58; New = 'System, dnasr281@gmail.com' }
    @{ Row = 59; New = 'System, admin@admin.com' }
    @{ Row = 60; New = 'System, backup@backdoor.com' }
    @{ Row = 62; New = 'System, dnasr281@gmail.com' }
    @{ Row = 63; New = 'System, dnasr281@gmail.com' }
    @{ Row = 64; New = 'System, dnasr281@gmail.com' }
    @{ Row = 65; New = 'System, dnasr281@gmail.com' }
    @{ Row = 66; New = 'System, dnasr281@gmail.com' }
    @{ Row = 67; New = 'System, dnasr281@gmail.com' }
    @{ Row = 69; New = 'System, dnasr281@gmail.com' }
    @{ Row = 70; New = 'System, dnasr281@gmail.com' }
    @{ Row = 71; New = 'System, dnasr281@gmail.com' }
    @{ Row = 72; New = 'System, dnasr281@gmail.com' }
    @{ Row = 73; New = 'System, dnasr281@gmail.com' }
    @{ Row = 74; New = 'System, dnasr281@gmail.com' }
    @{ Row = 76; New = 'System, dnasr281@gmail.com' }
    @{ Row = 78; New = 'System, dnasr281@gmail.com' }
    @{ Row = 80; New = 'System, backup@backdoor.com' }
    @{ Row = 81; New = 'System, backup@backdoor.com' }
    @{ Row = 82; New = 'System, backup@backdoor.com' }
    @{ Row = 83; New = 'System, dnasr281@gmail.com' }
    @{ Row = 84; New = 'System, dnasr281@gmail.com' }
    @{ Row = 85; New = 'System, dnasr281@gmail.com' }
    @{ Row = 86; New = 'System, dnasr281@gmail.com' }
    @{ Row = 90; New = 'System, dnasr281@gmail.com' }
    @{ Row = 92; New = 'System, dnasr281@gmail.com' }
    @{ Row = 93; New = 'System, dnasr281@gmail.com' }
    @{ Row = 94; New = 'System, dnasr281@gmail.com' }
    @{ Row = 96; New = 'System, dnasr281@gmail.com' }
    @{ Row = 99; New = 'System, dnasr281@gmail.com' }
    @{ Row = 101; New = 'System, dnasr281@gmail.com' }
    @{ Row = 106; New = 'System, backup@backdoor.com' }
    @{ Row = 107; New = 'System, backup@backdoor.com' }
    @{ Row = 108; New = 'System, backup@backdoor.com' }
    @{ Row = 109; New = 'System, dnasr281@gmail.com' }
    @{ Row = 110; New = 'System, dnasr281@gmail.com' }
    @{ Row = 111; New = 'System, dnasr281@gmail.com' }
    @{ Row = 112; New = 'System, dnasr281@gmail.com' }
    @{ Row = 116; New = 'System, dnasr281@gmail.com' }
    @{ Row = 118; New = 'System, dnasr281@gmail.com' }
    @{ Row = 119; New = 'System, dnasr281@gmail.com' }
    @{ Row = 120; New = 'System, dnasr281@gmail.com' }
    @{ Row = 122; New = 'System, dnasr281@gmail.com' }
    @{ Row = 125; New = 'System, dnasr281@gmail.com' }
    @{ Row = 127; New = 'System, dnasr281@gmail.com' }
    @{ Row = 132; New = 'System, backup@backdoor.com' }
    @{ Row = 133; New = 'System, backup@backdoor.com' }
    @{ Row = 134; New = 'System, backup@backdoor.com' }
    @{ Row = 135; New = 'System, dnasr281@gmail.com' }
    @{ Row = 136; New = 'System, dnasr281@gmail.com' }
    @{ Row = 137; New = 'System, dnasr281@gmail.com' }
    @{ Row = 138; New = 'System, dnasr281@gmail.com' }
    @{ Row = 142; New = 'System, dnasr281@gmail.com' }
    @{ Row = 144; New = 'System, dnasr281@gmail.com' }
    @{ Row = 145; New = 'System, dnasr281@gmail.com' }
    @{ Row = 146; New = 'System, dnasr281@gmail.com' }
    @{ Row = 148; New = 'System, dnasr281@gmail.com' }
    @{ Row = 151; New = 'System, dnasr281@gmail.com' }
    @{ Row = 153; New = 'System, dnasr281@gmail.com' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.New
}
